$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 118, shifting existing rows 118:164 down to 119:165
$ws.Rows.Item(118).Insert()

# Populate the new row 118 with the latest weekly data point
$ws.Cells.Item(118, 1).Value = 10
$ws.Cells.Item(118, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(118, 3).Value = "La Araucanía"
$ws.Cells.Item(118, 4).Value = 44704
$ws.Cells.Item(118, 4).NumberFormat = $ws.Cells.Item(119, 4).NumberFormat
$ws.Cells.Item(118, 5).Value = 9
$ws.Cells.Item(118, 6).Value = 100114007
$ws.Cells.Item(118, 7).Value = "Jengibre"
$ws.Cells.Item(118, 8).Value = "Sin especificar"
$ws.Cells.Item(118, 9).Value = "Primera"
$ws.Cells.Item(118, 10).Value = 40
$ws.Cells.Item(118, 11).Value = 20000
$ws.Cells.Item(118, 12).Value = 20000
$ws.Cells.Item(118, 13).Value = 20000
$ws.Cells.Item(118, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(118, 15).Value = "Perú"
$ws.Cells.Item(118, 16).Value = 1538
$ws.Cells.Item(118, 17).Value = 13
$ws.Cells.Item(118, 18).Value = "Hortaliza"
